$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.040.42'
$ws.Range("E2").Value = '  +6.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.669.01'
$ws.Range("E3").Value = '  +18.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.30'
$ws.Range("E5").Value = '  +7.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.87'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.664.81'
$ws.Range("E7").Value = '  +18.08%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +5.86%  '
$ws.Range("E10").Value = '  +8.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.69'
$ws.Range("E11").Value = '  +5.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("E12").Value = '  +7.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.51'
$ws.Range("E13").Value = '  +11.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000255'
$ws.Range("E14").Value = '  +6.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.280.24'
$ws.Range("E15").Value = '  +18.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '71.005.38'
$ws.Range("E16").Value = '  +6.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.663.16'
$ws.Range("E17").Value = '  +18.08%  '
$ws.Range("E18").Value = '  +2.01%  '
$ws.Range("E19").Value = '  +7.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '522.38'
$ws.Range("E20").Value = '  +8.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.95'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.32'
$ws.Range("E22").Value = '  +19.28%  '
$ws.Range("E23").Value = '  +7.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.55'
$ws.Range("E24").Value = '  +13.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.82'
$ws.Range("E25").Value = '  +6.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.53'
$ws.Range("E26").Value = '  +7.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.13'
$ws.Range("E27").Value = '  +9.80%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  +11.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.21'
$ws.Range("E30").Value = '  +3.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.92'
$ws.Range("E31").Value = '  +11.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.77'
$ws.Range("E32").Value = '  +13.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0000111'
$ws.Range("E33").Value = '  +17.82%  '
$ws.Range("E34").Value = '  +4.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.14'
$ws.Range("E36").Value = '  +9.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.353'
$ws.Range("E37").Value = '  +13.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.03'
$ws.Range("E38").Value = '  +9.27%  '
$ws.Range("E39").Value = '  +9.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("E40").Value = '  +6.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.27'
$ws.Range("E41").Value = '  +4.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '46.14'
$ws.Range("E42").Value = '  -4.72%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.13'
$ws.Range("E43").Value = '  +13.96%  '
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.84'
$ws.Range("E44").Value = '  +5.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.119.26'
$ws.Range("E45").Value = '  +11.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("E47").Value = '  +7.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.47'
$ws.Range("E48").Value = '  +11.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.78'
$ws.Range("E49").Value = '  +3.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("E51").Value = '  +10.99%  '
